$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Incidenten")
$ws4 = $wb.Worksheets.Item("Handelingen")

$ws4.Select()
$ws1.Select()
$ws1.Range("B10").Select()
$ws4.Select()
Write-Host "done"
